$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K58").Value = 1634355.84
$ws.Range("G59").Value = -885267.008
$ws.Range("K59").Value = -876074.048
$ws.Range("G60").Value = 725972.096
$ws.Range("G61").Value = -697479.872
$ws.Range("K61").Value = -763313.024
$ws.Range("G62").Value = 186089.04
$ws.Range("K62").Value = 228795.04
$ws.Range("G63").Value = -332008.992
$ws.Range("K63").Value = -343604.064
$ws.Range("O63").Value = -353134.016
$ws.Range("K64").Value = -303207.968
$ws.Range("G65").Value = -65637
$ws.Range("K65").Value = -72281.008
$ws.Range("O65").Value = -80414.992
$ws.Range("G66").Value = 31915
$ws.Range("G69").Value = 28491.992
$ws.Range("G73").Value = 7562.008
$ws.Range("K73").Value = -27774.992
$ws.Range("G74").Value = -36219.992
$ws.Range("K74").Value = 14457
$ws.Range("O74").Value = -672.024
$ws.Range("G79").Value = 35181.008
$ws.Range("K79").Value = -14500.984
$ws.Range("O79").Value = 50958.04
